# Merging of suites and updation of code
# Update the BOL sheet's "Way Bill #" test-data order IDs (column A,
# rows 2-4) with the newer order numbers from the merged test suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOL")

$ws.Range("A2").Value = "51500899"
$ws.Range("A3").Value = "51487922"
$ws.Range("A4").Value = "51487922"
